# "Function to send email with the grades created. Missing attaching file to the email"
#
# The author is preparing to email each student's grades. As part of that
# work, a consolidated "samples" worksheet is created at the end of the
# workbook holding a full copy of both class rosters (Matematicas and
# Desarrollo Web) side by side, while the original per-class sheets are
# trimmed back down to just a couple of sample rows.

$wb = $excel.ActiveWorkbook

$mathWs = $wb.Worksheets.Item("Matematicas")
$webWs  = $wb.Worksheets.Item("Desarrollo Web")

# ---------------------------------------------------------------------------
# 1. Add a new worksheet named "samples" as the last tab of the workbook.
# ---------------------------------------------------------------------------
$lastWs     = $wb.Worksheets.Item($wb.Worksheets.Count)
$samplesWs  = $wb.Worksheets.Add($null, $lastWs)
$samplesWs.Name = "samples"

# ---------------------------------------------------------------------------
# 2. Copy the full original rosters into "samples":
#      Desarrollo Web (A1:D11) -> samples A1:D11
#      Matematicas    (A1:D11) -> samples F1:I11
# ---------------------------------------------------------------------------
$webWs.Range("A1:D11").Copy()
$samplesWs.Range("A1:D11").PasteSpecial(-4163)

$mathWs.Range("A1:D11").Copy()
$samplesWs.Range("F1:I11").PasteSpecial(-4163)

$excel.CutCopyMode = 0
$samplesWs.Range("G18").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Trim "Desarrollo Web" back down to the header + first two students.
# ---------------------------------------------------------------------------
$webWs.Range("A4:D11").ClearContents()
$webWs.Range("A4:XFD11").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Trim "Matematicas" back down to the header + first two students, and
#    reserve F19:G19 (left blank, ready for the attachment helper cells).
#    Activate this sheet last so it ends up the selected tab.
# ---------------------------------------------------------------------------
$mathWs.Range("A4:D11").ClearContents()
$mathWs.Range("F19:G19").Font.Bold = $false

[void]$mathWs.Activate()
$mathWs.Range("F19:G19").Select() | Out-Null
